# Auto-generated edit script: update crypto price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.111.68'
$ws.Range('E2').Value = '  +1.51%  '
$ws.Range('D3').Value = '3.178.41'
$ws.Range('E3').Value = '  +3.61%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.57'
$ws.Range('E5').Value = '  +3.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.96'
$ws.Range('E6').Value = '  +4.44%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '3.178.02'
$ws.Range('E8').Value = '  +3.69%  '
$ws.Range('E10').Value = '  +5.05%  '
$ws.Range('E11').Value = '  +0.78%  '
$ws.Range('E12').Value = '  +2.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000274'
$ws.Range('E13').Value = '  +18.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '38.08'
$ws.Range('E14').Value = '  +6.41%  '
$ws.Range('D15').Value = '3.699.29'
$ws.Range('E15').Value = '  +3.84%  '
$ws.Range('D16').Value = '65.198.31'
$ws.Range('E16').Value = '  +1.57%  '
$ws.Range('D17').Value = '3.181.65'
$ws.Range('E17').Value = '  +3.65%  '
$ws.Range('E18').Value = '  +5.63%  '
$ws.Range('E19').Value = '  +1.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '514.03'
$ws.Range('E20').Value = '  +7.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.91'
$ws.Range('E21').Value = '  +6.64%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.734'
$ws.Range('E22').Value = '  +7.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '15.36'
$ws.Range('E23').Value = '  +6.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.83'
$ws.Range('E24').Value = '  +3.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.47'
$ws.Range('E25').Value = '  +3.46%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.05'
$ws.Range('E27').Value = '  +10.63%  '
$ws.Range('E28').Value = '  +3.97%  '
$ws.Range('E29').Value = '  +6.73%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '28.08'
$ws.Range('E30').Value = '  +6.42%  '
$ws.Range('E31').Value = '  +13.21%  '
$ws.Range('E32').Value = '  +0.13%  '
$ws.Range('E33').Value = '  +5.84%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.25'
$ws.Range('E34').Value = '  +7.07%  '
$ws.Range('E35').Value = '  +5.66%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '55.72'
$ws.Range('E36').Value = '  +1.40%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0902'
$ws.Range('E37').Value = '  +10.19%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '478.40'
$ws.Range('E38').Value = '  +5.94%  '
$ws.Range('E39').Value = '  +10.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0423'
$ws.Range('E40').Value = '  +2.47%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.66'
$ws.Range('E41').Value = '  +4.27%  '
$ws.Range('D42').Value = '3.069.30'
$ws.Range('E42').Value = '  +1.23%  '
$ws.Range('E43').Value = '  +3.39%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.288'
$ws.Range('E44').Value = '  +7.23%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.40'
$ws.Range('E45').Value = '  +6.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '29.23'
$ws.Range('E46').Value = '  +4.46%  '
$ws.Range('E47').Value = '  +15.96%  '
$ws.Range('E48').Value = '  -0.07%  '
$ws.Range('E49').Value = '  +2.21%  '
$ws.Range('E50').Value = '  +8.51%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '120.62'
$ws.Range('E51').Value = '  +1.76%  '
